$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 539, shifting existing rows 539:613 down to 540:614
$ws.Rows("539:539").Insert()

# Populate the newly inserted row 539 with the new record
$ws.Cells.Item(539, 1).Value = 4
$ws.Cells.Item(539, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(539, 3).Value = "Los Lagos"
$ws.Cells.Item(539, 4).Value = 44984
$ws.Cells.Item(539, 5).Value = 10
$ws.Cells.Item(539, 6).Value = 100114001
$ws.Cells.Item(539, 7).Value = "Papa"
$ws.Cells.Item(539, 8).Value = "Red Lady"
$ws.Cells.Item(539, 9).Value = "1a (cosecha)"
$ws.Cells.Item(539, 10).Value = 250
$ws.Cells.Item(539, 11).Value = 12000
$ws.Cells.Item(539, 12).Value = 13000
$ws.Cells.Item(539, 13).Value = 12600
$ws.Cells.Item(539, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(539, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(539, 16).Value = 504
$ws.Cells.Item(539, 17).Value = 25
$ws.Cells.Item(539, 18).Value = "Hortaliza"

# Keep the same date number format as the other date cells in column D
$ws.Cells.Item(539, 4).NumberFormat = $ws.Cells.Item(540, 4).NumberFormat
